$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1569.3182
$ws.Range("J58").Value = 1728.5714
$ws.Range("L58").Value = 5185.7142
$ws.Range("N58").Value = -5485.7142

$ws.Range("H80").Value = 1656.5
$ws.Range("I80").Value = 1496.2
$ws.Range("J80").Value = 1745.5555
$ws.Range("K80").Value = 4488.6
$ws.Range("L80").Value = 5236.666499999999
$ws.Range("M80").Value = -3490.6
$ws.Range("N80").Value = -7232.666499999999

$ws.Range("H83").Value = 1656.5
$ws.Range("I83").Value = 1496.2
$ws.Range("J83").Value = 1745.5555
$ws.Range("K83").Value = 13465.8
$ws.Range("L83").Value = 15709.9995
$ws.Range("M83").Value = -8473.800000000001
$ws.Range("N83").Value = -25693.9995

$ws.Range("H106").Value = 671503.5600000001
$ws.Range("I106").Value = 1433049.6
$ws.Range("J106").Value = 5150.75
$ws.Range("K106").Value = 1433049.6
$ws.Range("L106").Value = 5150.75
$ws.Range("M106").Value = -1432418.6
$ws.Range("N106").Value = -6412.75

$ws.Range("H132").Value = 4978.871
$ws.Range("I132").Value = 1938.409
$ws.Range("J132").Value = 12411.111
$ws.Range("K132").Value = 5815.227000000001
$ws.Range("L132").Value = 37233.333
$ws.Range("M132").Value = -3285.227000000001
$ws.Range("N132").Value = -42293.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1487.3529
$ws.Range("I2").Value = 1351.0526
$ws.Range("J2").Value = 1660
$ws.Range("K2").Value = 1351.0526
$ws.Range("L2").Value = 1660
$ws.Range("M2").Value = -1238.0526
$ws.Range("N2").Value = -1886

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H61").Value = 10347.647
$ws.Range("I61").Value = 1560.6666
$ws.Range("K61").Value = 1560.6666
$ws.Range("M61").Value = -1348.6666

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H110").Value = 949.2222
$ws.Range("I110").Value = 840
$ws.Range("J110").Value = 1036.6
$ws.Range("K110").Value = 840
$ws.Range("L110").Value = 1036.6
$ws.Range("M110").Value = 1205
$ws.Range("N110").Value = -5126.6

$ws.Range("H116").Value = 1487.3529
$ws.Range("I116").Value = 1351.0526
$ws.Range("J116").Value = 1660
$ws.Range("K116").Value = 1351.0526
$ws.Range("L116").Value = 1660
$ws.Range("M116").Value = 942.9474
$ws.Range("N116").Value = -6248

$ws.Range("H136").Value = 10347.647
$ws.Range("I136").Value = 1560.6666
$ws.Range("K136").Value = 4681.9998
$ws.Range("M136").Value = -2131.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1487.3529
$ws.Range("I3").Value = 1351.0526
$ws.Range("J3").Value = 1660
$ws.Range("K3").Value = 1351.0526
$ws.Range("L3").Value = 1660
$ws.Range("M3").Value = -1237.0526
$ws.Range("N3").Value = -1888

$ws.Range("H105").Value = 2309.375
$ws.Range("I105").Value = 2245.8333
$ws.Range("K105").Value = 2245.8333
$ws.Range("M105").Value = -498.8332999999998

$ws.Range("H107").Value = 1148.2727
$ws.Range("I107").Value = 804.4286
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 804.4286
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 1115.5714
$ws.Range("N107").Value = -5590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2046.1167
$ws.Range("I31").Value = 1288.9667
$ws.Range("J31").Value = 2803.2666
$ws.Range("K31").Value = 1288.9667
$ws.Range("L31").Value = 2803.2666
$ws.Range("M31").Value = -993.9666999999999
$ws.Range("N31").Value = -3393.2666

$ws.Range("H34").Value = 2046.1167
$ws.Range("I34").Value = 1288.9667
$ws.Range("J34").Value = 2803.2666
$ws.Range("K34").Value = 1288.9667
$ws.Range("L34").Value = 2803.2666
$ws.Range("M34").Value = -1086.9667
$ws.Range("N34").Value = -3207.2666

$ws.Range("H132").Value = 1196343.9
$ws.Range("I132").Value = 1367.7916
$ws.Range("J132").Value = 5293404.5
$ws.Range("K132").Value = 4103.3748
$ws.Range("L132").Value = 15880213.5
$ws.Range("M132").Value = -1573.3748
$ws.Range("N132").Value = -15885273.5

$ws.Range("H141").Value = 63485.355
$ws.Range("J141").Value = 63485.355
$ws.Range("L141").Value = 63485.355
$ws.Range("N141").Value = -73845.35500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 522.2558
$ws.Range("I5").Value = 259.125
$ws.Range("K5").Value = 777.375
$ws.Range("M5").Value = -665.375

$ws.Range("H69").Value = 2898.1
$ws.Range("J69").Value = 3164.5557
$ws.Range("L69").Value = 9493.667099999999
$ws.Range("N69").Value = -11115.6671

$ws.Range("H72").Value = 2898.1
$ws.Range("J72").Value = 3164.5557
$ws.Range("L72").Value = 28481.0013
$ws.Range("N72").Value = -36593.0013

$ws.Range("H131").Value = 1308.3
$ws.Range("I131").Value = 577
$ws.Range("J131").Value = 1330.9175
$ws.Range("K131").Value = 1731
$ws.Range("L131").Value = 3992.7525
$ws.Range("M131").Value = 3309
$ws.Range("N131").Value = -14072.7525

$ws.Range("H135").Value = 522.2558
$ws.Range("I135").Value = 259.125
$ws.Range("K135").Value = 2332.125
$ws.Range("M135").Value = 202.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 18251
$ws.Range("J93").Value = 18251
$ws.Range("L93").Value = 18251
$ws.Range("N93").Value = -21995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 740
$ws.Range("I82").Value = 740
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 740
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -379
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 740
$ws.Range("I85").Value = 740
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 740
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 508
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5054.175
$ws.Range("I136").Value = 2698.3462
$ws.Range("J136").Value = 9429.286
$ws.Range("K136").Value = 8095.0386
$ws.Range("L136").Value = 28287.858
$ws.Range("M136").Value = -5545.0386
$ws.Range("N136").Value = -33387.858
